$d = $word.ActiveDocument

$replacements = @(
    @("94÷5=", "40÷5="),
    @("60÷7=", "43÷7="),
    @("58÷2=", "83÷8="),
    @("37÷7=", "20÷4="),
    @("16÷7=", "88÷8="),
    @("70÷4=", "33÷6="),
    @("56÷5=", "83÷6="),
    @("41÷4=", "11÷6="),
    @("13÷7=", "61÷6="),
    @("13÷6=", "56÷8="),
    @("17÷6=", "83÷9="),
    @("83÷7=", "34÷5="),
    @("50÷7=", "48÷2="),
    @("58÷8=", "29÷9="),
    @("46÷8=", "15÷5="),
    @("25÷4=", "91÷3="),
    @("76÷2=", "29÷3="),
    @("45÷9=", "94÷9="),
    @("60÷3=", "56÷6="),
    @("24÷8=", "48÷9="),
    @("83÷5=", "95÷4="),
    @("25÷3=", "57÷8="),
    @("39÷2=", "73÷6="),
    @("63÷2=", "92÷2="),
    @("31÷3=", "25÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
